# Fill in the BMI Value (column C) and BMI Category (column D) for each
# person's row, computed from their Height (A) and Weight (B) already
# present in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bmiValues = @{
    2 = 179.59
    3 = 210.65
    4 = 918.27
    5 = 208.28
    6 = 340.9
    7 = 71.02
    8 = 303.02999999999997
}

foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = $bmiValues[$row]
    $ws.Cells.Item($row, 4).Value = "obese category"
}

# Mirror the author's final cursor position (just below the last data row).
$ws.Range("C9").Select()
